# cambios en todo 26-7
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PreparacionDatos")

# ---------------------------------------------------------------------
# 1) Remove the pre-existing hyperlink (it sat on I2, which is being
#    repurposed); we'll re-add equivalent hyperlinks at their new homes.
# ---------------------------------------------------------------------
foreach ($h in $ws.Hyperlinks) { $h.Delete() }

# ---------------------------------------------------------------------
# 2) The old data row (row 2) moves down to row 4 unchanged (values +
#    formatting), to make room for two freshly entered rows above it.
# ---------------------------------------------------------------------
foreach ($col in @("A","B","C","D","E","F","G","H","I")) {
    $src = $ws.Range($col + "2")
    $dst = $ws.Range($col + "4")
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $dst.Value2 = $src.Value2
}
$ws.Application.CutCopyMode = $false

# Correct the document number text on the relocated row (543416748666 -> 543416876777)
$ws.Range("C4").Value2 = "543416876777"

# ---------------------------------------------------------------------
# 3) Carry over the date / hyperlink-email formatting that used to live
#    in H2/I2 into the new F/G columns (Fecha Nacimiento / Email) for
#    the two freshly inserted rows, before we overwrite H2/I2 below.
# ---------------------------------------------------------------------
$ws.Range("H2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F3").PasteSpecial(-4122)

$ws.Range("I2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G3").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4) Header row restyle: reuse the header format from C1 everywhere it
#    is needed (B1 now "N.º Documento"; J1:N1 are brand-new headers).
# ---------------------------------------------------------------------
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("L1").PasteSpecial(-4122)
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("N1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("B1").Value2 = "N.º Documento"
$ws.Range("C1").Value2 = "Nombre"
$ws.Range("D1").Value2 = "Apellido"
$ws.Range("E1").Value2 = "Sexo"
$ws.Range("F1").Value2 = "Fecha Nacimiento"
$ws.Range("G1").Value2 = "Email"
$ws.Range("H1").Value2 = "Plan"
$ws.Range("I1").Value2 = "Provincia"
$ws.Range("J1").Value2 = "Localidad"
$ws.Range("K1").Value2 = "Linea"
$ws.Range("L1").Value2 = "ICCID"
$ws.Range("M1").Value2 = "IMSI"
$ws.Range("N1").Value2 = "KI"

# ---------------------------------------------------------------------
# 5) New row 2 ("Malan Fareto" line) + row 3 ("Natias Mazano" line).
# ---------------------------------------------------------------------
$ws.Range("K2").NumberFormat = "@"
$ws.Range("L2").NumberFormat = "@"
$ws.Range("M2").NumberFormat = "@"
$ws.Range("N2").NumberFormat = "@"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("K3").NumberFormat = "@"
$ws.Range("L3").NumberFormat = "@"
$ws.Range("M3").NumberFormat = "@"
$ws.Range("N3").NumberFormat = "@"
$ws.Range("B3").NumberFormat = "@"

$ws.Range("A2").Value2 = "Alta Linea"
$ws.Range("B2").Value2 = "59885659"
$ws.Range("C2").Value2 = "Malan"
$ws.Range("D2").Value2 = "Fareto"
$ws.Range("E2").Value2 = "Masculino"
$ws.Range("F2").Value2 = "08/08/1992"
$ws.Range("G2").Value2 = "malannominacion@gmail.com"
$ws.Range("H2").Value2 = "Plan prepago nacional"
$ws.Range("I2").Value2 = "Buenos Aires"
$ws.Range("J2").Value2 = "vicente lopez"
$ws.Range("K2").Value2 = "52648975215"
$ws.Range("L2").Value2 = "15425684"
$ws.Range("M2").Value2 = "25489645"
$ws.Range("N2").Value2 = "14756841"

$ws.Range("A3").Value2 = "Alta Linea"
$ws.Range("B3").Value2 = "59885660"
$ws.Range("C3").Value2 = "Natias"
$ws.Range("D3").Value2 = "Mazano"
$ws.Range("E3").Value2 = "Masculino"
$ws.Range("F3").Value2 = "08/08/1992"
$ws.Range("G3").Value2 = "malannominacion@gmail.com"
$ws.Range("H3").Value2 = "Plan prepago nacional"
$ws.Range("I3").Value2 = "Buenos Aires"
$ws.Range("J3").Value2 = "vicente lopez"
$ws.Range("K3").Value2 = "52648975215"
$ws.Range("L3").Value2 = "15425684"
$ws.Range("M3").Value2 = "25489645"
$ws.Range("N3").Value2 = "14756841"

# ---------------------------------------------------------------------
# 6) Hyperlinks for the three mailto cells.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:malannominacion@gmail.com")
$ws.Hyperlinks.Add($ws.Range("I4"), "mailto:malannominacion@gmail.com")
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:malannominacion@gmail.com")

# ---------------------------------------------------------------------
# 7) Column widths for the two newly introduced columns.
# ---------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 21.28515625
$ws.Columns("G").ColumnWidth = 26

# ---------------------------------------------------------------------
# 8) Sheet view: selection moves, this sheet is no longer the active tab.
# ---------------------------------------------------------------------
$ws.Range("D9").Select()

# ---------------------------------------------------------------------
# 9) "OM" becomes the active sheet / tab.
# ---------------------------------------------------------------------
$omWs = $wb.Worksheets.Item("OM")
$omWs.Activate()
$omWs.Range("H58").Select()
